$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "tbCuentaBancaria" block (rows 1-5) down to rows 8-12,
# opening space at the top for the new "tbBanco" block.
$ws.Rows("1:7").Insert()

# The old row 3 ("NombreBanco"/50) becomes part of the new "tbBanco" block
# instead, so remove it from the shifted-down block (this also pulls rows
# 11-12 up to 10-11).
$ws.Rows("10:10").Delete()

# In the shifted "tbCuentaBancaria" block, make room for the new "Tipo"
# column (B) by moving the PK/Identity/value cells one column to the right,
# and mark the (empty) B8 cell bold to match the header row A8.
$ws.Range("B8").Font.Bold = $true

$ws.Range("D9").Value = "Identity"
$ws.Range("C9").Value = "PK"
$ws.Range("B9").ClearContents()

$ws.Range("C10").Value = 20
$ws.Range("B10").ClearContents()

# ---- New block: tbBanco table definition (rows 1-3) ----
$ws.Range("A1").Value = "tbBanco"
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1").Font.Bold = $true

$ws.Range("A2").Value = "IdBanco"
$ws.Range("B2").Value = "int"
$ws.Range("C2").Value = "PK"
$ws.Range("D2").Value = "Identity"

$ws.Range("A3").Value = "NombreBanco"
$ws.Range("B3").Value = "string"
$ws.Range("C3").Value = 50

# Widen the new "Tipo" column (B) to match column A's width.
$ws.Columns("B").ColumnWidth = 15.877604166666666

# Match the saved selection/active cell from the target file.
$ws.Range("A10").Select()
